# Refresh the cryptos price list (columns D/E) with the latest scraped values,
# and fix the NEARProtocol/Quant row ordering (rows 47-48) to match the new data.
# Numeric-looking text (prices, e.g. "1.013") is written with a leading apostrophe
# so Excel keeps storing it as literal text instead of converting it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'27.513.03"
$ws.Range("E2").Value = "`'  +2.12%  "

$ws.Range("D3").Value = "`'1.868.36"
$ws.Range("E3").Value = "`'  +1.17%  "

$ws.Range("D4").Value = "`'1.013"
$ws.Range("E4").Value = "`'  +0.24%  "

$ws.Range("D5").Value = "`'311.94"
$ws.Range("E5").Value = "`'  +0.62%  "

$ws.Range("E6").Value = "`'  +0.32%  "

$ws.Range("D7").Value = "`'0.4780"
$ws.Range("E7").Value = "`'  +0.00%  "

$ws.Range("D8").Value = "`'0.3735"
$ws.Range("E8").Value = "`'  +1.76%  "

$ws.Range("D9").Value = "`'0.07314"
$ws.Range("E9").Value = "`'  +1.21%  "

$ws.Range("D10").Value = "`'0.9359"
$ws.Range("E10").Value = "`'  +0.94%  "

$ws.Range("D11").Value = "`'20.69"
$ws.Range("E11").Value = "`'  +5.06%  "

$ws.Range("D12").Value = "`'0.07835"
$ws.Range("E12").Value = "`'  +1.61%  "

$ws.Range("D13").Value = "`'1.859.46"
$ws.Range("E13").Value = "`'  +1.18%  "

$ws.Range("D14").Value = "`'5.441"
$ws.Range("E14").Value = "`'  +2.34%  "

$ws.Range("D15").Value = "`'6.554"
$ws.Range("E15").Value = "`'  +2.18%  "

$ws.Range("D16").Value = "`'90.17"
$ws.Range("E16").Value = "`'  +1.40%  "

$ws.Range("E17").Value = "`'  +0.18%  "

$ws.Range("D18").Value = "`'0.000008903"
$ws.Range("E18").Value = "`'  +2.98%  "

$ws.Range("E19").Value = "`'  +0.18%  "

$ws.Range("D20").Value = "`'27.533.46"
$ws.Range("E20").Value = "`'  +2.09%  "

$ws.Range("D21").Value = "`'14.63"
$ws.Range("E21").Value = "`'  +0.53%  "

$ws.Range("D22").Value = "`'5.117"
$ws.Range("E22").Value = "`'  +1.13%  "

$ws.Range("E23").Value = "`'  +0.44%  "

$ws.Range("E24").Value = "`'  +1.06%  "

$ws.Range("D25").Value = "`'154.41"
$ws.Range("E25").Value = "`'  +1.30%  "

$ws.Range("D26").Value = "`'18.46"
$ws.Range("E26").Value = "`'  +1.59%  "

$ws.Range("D27").Value = "`'2.021"
$ws.Range("E27").Value = "`'  +1.30%  "

$ws.Range("D28").Value = "`'115.65"
$ws.Range("E28").Value = "`'  +1.34%  "

$ws.Range("D29").Value = "`'4.987"

$ws.Range("D30").Value = "`'0.08912"
$ws.Range("E30").Value = "`'  +0.34%  "

$ws.Range("D31").Value = "`'3.337"
$ws.Range("E31").Value = "`'  +0.56%  "

$ws.Range("E32").Value = "`'  +3.58%  "

$ws.Range("D33").Value = "`'0.7582"
$ws.Range("E33").Value = "`'  +1.85%  "

$ws.Range("D34").Value = "`'4.613"

$ws.Range("D35").Value = "`'2.736"
$ws.Range("E35").Value = "`'  +0.72%  "

$ws.Range("D36").Value = "`'0.02043"
$ws.Range("E36").Value = "`'  +4.33%  "

$ws.Range("D37").Value = "`'1.119"
$ws.Range("E37").Value = "`'  -0.92%  "

$ws.Range("D38").Value = "`'2.999"
$ws.Range("E38").Value = "`'  +0.56%  "

$ws.Range("D39").Value = "`'0.05270"
$ws.Range("E39").Value = "`'  +0.10%  "

$ws.Range("D40").Value = "`'0.5310"
$ws.Range("E40").Value = "`'  +2.14%  "

$ws.Range("D41").Value = "`'7.071"
$ws.Range("E41").Value = "`'  +0.94%  "

$ws.Range("D42").Value = "`'0.1524"
$ws.Range("E42").Value = "`'  +0.82%  "

$ws.Range("D43").Value = "`'8.478"
$ws.Range("E43").Value = "`'  +3.33%  "

$ws.Range("D44").Value = "`'10.61"
$ws.Range("E44").Value = "`'  +0.52%  "

$ws.Range("D45").Value = "`'0.4801"
$ws.Range("E45").Value = "`'  +1.59%  "

$ws.Range("D46").Value = "`'1.014"
$ws.Range("E46").Value = "`'  +0.27%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "`'1.655"
$ws.Range("E47").Value = "`'  +3.13%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "`'102.78"
$ws.Range("E48").Value = "`'  +1.17%  "

$ws.Range("D49").Value = "`'67.37"
$ws.Range("E49").Value = "`'  +2.20%  "

$ws.Range("D50").Value = "`'0.06082"
$ws.Range("E50").Value = "`'  +0.92%  "

$ws.Range("D51").Value = "`'0.9196"
$ws.Range("E51").Value = "`'  +3.69%  "
